# cv121101a.xlsx — "correção nos dados e inicio da analise PNAD 2009"
#
# The original sheet had two section-header rows ("situação do domicílio" at
# row 5 and "grandes regiões e unidades da federação" at row 8) that carried
# a label but no data; the numbers for "urbana"/"rural"/each
# region/state had all been entered one row too low as a result. The fix
# removes those two spacer/header rows entirely so every label lines up with
# its own numbers, and relabels the "unnamed: 1_level_1" sub-header as
# "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "situação do domicílio" header row (row 5) — this shifts
# "urbana"/"rural" (and everything below) up by one, so they now carry the
# values that used to sit one row further down.
$ws.Rows.Item(5).Delete()

# Remove the "grandes regiões e unidades da federação" header row (now at
# row 7 after the previous deletion) — shifts "norte" and every
# region/state up by one more row.
$ws.Rows.Item(7).Delete()

# Relabel the second-header cell from "unnamed: 1_level_1" to "total".
$ws.Range("B2").Value = "total"
